$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 160 (this shifts the existing rows 160-260 down to 162-262,
# matching the target dimension A1:R262)
$ws.Rows("160:161").Insert()

# ---- Row 160 (new data) ----
$ws.Range("A160").Value = 9
$ws.Range("B160").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C160").Value = "Metropolitana"
$ws.Range("D160").Value = 44488
$ws.Range("E160").Value = 13
$ws.Range("F160").Value = 100112012
$ws.Range("G160").Value = "Espinaca"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 250
$ws.Range("K160").Value = 5000
$ws.Range("L160").Value = 6000
$ws.Range("M160").Value = 5500
$ws.Range("N160").Value = '$/cuna 10 kilos'
$ws.Range("O160").Value = "Provincia de Chacabuco"
$ws.Range("P160").Value = 550
$ws.Range("Q160").Value = 10
$ws.Range("R160").Value = "Hortaliza"

# ---- Row 161 (new data) ----
$ws.Range("A161").Value = 9
$ws.Range("B161").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C161").Value = "Metropolitana"
$ws.Range("D161").Value = 44488
$ws.Range("E161").Value = 13
$ws.Range("F161").Value = 100112012
$ws.Range("G161").Value = "Espinaca"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 97
$ws.Range("K161").Value = 4000
$ws.Range("L161").Value = 4000
$ws.Range("M161").Value = 4000
$ws.Range("N161").Value = '$/cuna 10 kilos'
$ws.Range("O161").Value = "Provincia de Chacabuco"
$ws.Range("P161").Value = 400
$ws.Range("Q161").Value = 10
$ws.Range("R161").Value = "Hortaliza"

Write-Host "Edit applied"
